# "db name change updates" - rename header cell N1 from
# cumulative_floor_surface_m2 -> address_floor_surface_m2, and update the
# sheet's active selection to match the saved view (O8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("buildings")

$ws.Range("N1").Value = "address_floor_surface_m2"

$ws.Activate()
$ws.Range("O8").Select()
